$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts Part-color/Technical-qc/Description right)
$ws.Columns("B:B").Insert()

# The freshly inserted column takes on the width of the column to its left (A)
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Set the new header value in the freshly inserted column B
$ws.Range("B1").Value = "Sp-Category"

# Match the active cell/selection shown in the target sheet
$ws.Range("B1").Select() | Out-Null
